$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label touch-ups (row 1) ---
$ws.Range("D1").Value = "SHARKweb english"
$ws.Range("E1").Value = "SHARKweb swedish"

# --- dataset / datatype_code rename (row 10) ---
$ws.Range("B10").Value = "datatype_code"

# --- Drop the "(TEST)" suffix from the Month export labels (row 14) ---
$ws.Range("G14").Value = "Month"
$ws.Range("H14").Value = "Månad"

# --- variable column_type renames ---
$ws.Range("B24").Value = "scientific_name"
$ws.Range("B34").Value = "dyntaxa_id"
$ws.Range("B42").Value = "reported_scientific_name"

# --- Selection / scroll position the author left the sheet in ---
$ws.Range("B24").Select()
